$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 109
$ws1.Range("F3").Value = 180
$ws1.Range("F5").Value = 199
$ws1.Range("F7").Value = 1156
$ws1.Range("F8").Value = 388
$ws1.Range("F9").Value = 197
$ws1.Range("F10").Value = 52
$ws1.Range("F13").Value = 400
$ws1.Range("F16").Value = 725
$ws1.Range("F19").Value = 1011
$ws1.Range("F20").Value = 462
$ws1.Range("F23").Value = 381
$ws1.Range("F25").Value = 42

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 366
$ws2.Range("F7").Value = 285

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 109
$ws4.Range("F5").Value = 180
$ws4.Range("F7").Value = 199
$ws4.Range("F9").Value = 1156
$ws4.Range("F10").Value = 388
$ws4.Range("F11").Value = 197
$ws4.Range("F13").Value = 52
$ws4.Range("F14").Value = 366
$ws4.Range("F19").Value = 285
$ws4.Range("F20").Value = 400
$ws4.Range("F23").Value = 725
$ws4.Range("F26").Value = 1011
$ws4.Range("F27").Value = 462
$ws4.Range("F32").Value = 381
$ws4.Range("F36").Value = 42
